$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$ws1.Name = "Test_Yearly"
$ws2.Name = "Test_Weekly"

$ws1.Select()
